# "added 4wk low sales check"
# Refresh the forecast numbers on the "Forecast Comparison" sheet (MyForecast,
# Inventory Coverage, Seasonality Index) and roll the new totals up into the
# "Summary" sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison: MyForecast (D), Inventory Coverage (H), Seasonality Index (L) ---

$ws1.Range("D2").Value = 63
$ws1.Range("H2").Value = 24.94
$ws1.Range("L2").Value = 1.17

$ws1.Range("D3").Value = 59
$ws1.Range("H3").Value = 25.64
$ws1.Range("L3").Value = 1.06

$ws1.Range("D4").Value = 57
$ws1.Range("H4").Value = 25.23
$ws1.Range("L4").Value = 1

$ws1.Range("D5").Value = 57
$ws1.Range("H5").Value = 24.57
$ws1.Range("L5").Value = 0.91

$ws1.Range("D6").Value = 55
$ws1.Range("H6").Value = 24.25
$ws1.Range("L6").Value = 0.9399999999999999

$ws1.Range("H7").Value = 23.77
$ws1.Range("L7").Value = 0.92

$ws1.Range("D8").Value = 52
$ws1.Range("H8").Value = 23.73
$ws1.Range("L8").Value = 1.15

$ws1.Range("H9").Value = 23
$ws1.Range("L9").Value = 0.88

$ws1.Range("H10").Value = 22.98
$ws1.Range("L10").Value = 0.93

$ws1.Range("D11").Value = 47
$ws1.Range("H11").Value = 22.62
$ws1.Range("L11").Value = 1.19

$ws1.Range("H12").Value = 21.71
$ws1.Range("L12").Value = 1

$ws1.Range("D13").Value = 46
$ws1.Range("H13").Value = 21.34
$ws1.Range("L13").Value = 0.8100000000000001

$ws1.Range("D14").Value = 45
$ws1.Range("H14").Value = 20.88
$ws1.Range("L14").Value = 1.15

$ws1.Range("D15").Value = 43
$ws1.Range("H15").Value = 20.71
$ws1.Range("L15").Value = 0.9

$ws1.Range("D16").Value = 42
$ws1.Range("H16").Value = 20.18
$ws1.Range("L16").Value = 0.83

$ws1.Range("D17").Value = 41
$ws1.Range("H17").Value = 19.55
$ws1.Range("L17").Value = 0.84

# --- Summary sheet: roll up totals (kept as text, matching the sheet's existing layout) ---

$summaryUpdates = @{
    "B9"  = "813"
    "B10" = "450"
    "B11" = "237"
    "B12" = "63"
    "B14" = "41"
}

foreach ($addr in $summaryUpdates.Keys) {
    $cell = $ws2.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $summaryUpdates[$addr]
}
